$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) First paragraph: "This is a Microsoft word document." gets two
#    trailing spaces, followed by a red (C00000) parenthetical note,
#    split across three runs exactly as the source edit produced them.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$body = $d.Range($p1.Range.Start, $p1.Range.End - 1)

$body.InsertAfter("  ")
$body.Collapse(0)

$enDash = [char]0x2013

$body.InsertAfter("(This is a change " + $enDash + " Ve")
$body.Font.Color = 192
$body.Collapse(0)

$body.InsertAfter("rsion for branch alternate")
$body.Font.Color = 192
$body.Collapse(0)

$body.InsertAfter(")")
$body.Font.Color = 192
$body.Collapse(0)

# ---------------------------------------------------------------------
# 2) Append a new, blank, shaded paragraph after the final paragraph of
#    the speech excerpt.
# ---------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$lastRange = $lastPara.Range

$lastRange.Find.Execute("free at last.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "free at last.^p", 2)

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.Style = $d.Styles("Normal")
$newPara.Format.Shading.Texture = 0
$newPara.Format.Shading.ForegroundPatternColor = -16777216
$newPara.Format.Shading.BackgroundPatternColor = 16382457
